$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 109.25
$ws.Range("I11").Value = 109.25
$ws.Range("K11").Value = 109.25
$ws.Range("M11").Value = 30.75

$ws.Range("H70").Value = 2029.7778
$ws.Range("J70").Value = 2129.5
$ws.Range("L70").Value = 6388.5
$ws.Range("N70").Value = -6928.5

$ws.Range("H73").Value = 2029.7778
$ws.Range("J73").Value = 2129.5
$ws.Range("L73").Value = 6388.5
$ws.Range("N73").Value = -8260.5

$ws.Range("H106").Value = 4373.5713
$ws.Range("I106").Value = 4373.5713
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4373.5713
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -3742.5713
$ws.Range("N106").ClearContents()

$ws.Range("H137").Value = 11948.75
$ws.Range("I137").Value = 1374.3846
$ws.Range("J137").Value = 31586.857
$ws.Range("K137").Value = 4123.1538
$ws.Range("L137").Value = 94760.571
$ws.Range("M137").Value = -1573.1538
$ws.Range("N137").Value = -99860.571

$ws.Range("H138").Value = 3866.644
$ws.Range("J138").Value = 4029.551
$ws.Range("L138").Value = 12088.653
$ws.Range("N138").Value = -22368.653

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2235011.5
$ws.Range("I32").Value = 2749063
$ws.Range("K32").Value = 2749063
$ws.Range("M32").Value = -2748776

$ws.Range("H61").Value = 51130.56
$ws.Range("I61").Value = 2221.919
$ws.Range("J61").Value = 352733.84
$ws.Range("K61").Value = 2221.919
$ws.Range("L61").Value = 352733.84
$ws.Range("M61").Value = -2009.919
$ws.Range("N61").Value = -353157.84

$ws.Range("H74").Value = 20546.857
$ws.Range("I74").Value = 1468.0588
$ws.Range("J74").Value = 101631.75
$ws.Range("K74").Value = 1468.0588
$ws.Range("L74").Value = 101631.75
$ws.Range("M74").Value = -594.0588
$ws.Range("N74").Value = -103379.75

$ws.Range("H77").Value = 20546.857
$ws.Range("I77").Value = 1468.0588
$ws.Range("J77").Value = 101631.75
$ws.Range("K77").Value = 7340.294
$ws.Range("L77").Value = 508158.75
$ws.Range("M77").Value = -2972.294
$ws.Range("N77").Value = -516894.75

$ws.Range("H132").Value = 2568386.8
$ws.Range("I132").Value = 1737.5555
$ws.Range("K132").Value = 5212.666499999999
$ws.Range("M132").Value = -2682.666499999999

$ws.Range("H136").Value = 51130.56
$ws.Range("I136").Value = 2221.919
$ws.Range("J136").Value = 352733.84
$ws.Range("K136").Value = 6665.757
$ws.Range("L136").Value = 1058201.52
$ws.Range("M136").Value = -4115.757
$ws.Range("N136").Value = -1063301.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 17905.615
$ws.Range("I99").Value = 18191.525
$ws.Range("K99").Value = 18191.525
$ws.Range("M99").Value = -16693.525

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 37915.668
$ws.Range("J50").Value = 37915.668
$ws.Range("L50").Value = 37915.668
$ws.Range("N50").Value = -39165.668

$ws.Range("H51").Value = 26665.5
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 29998.6
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 29998.6
$ws.Range("M51").Value = -9264
$ws.Range("N51").Value = -31470.6

$ws.Range("H60").Value = 19061.125
$ws.Range("I60").Value = 6498.4
$ws.Range("J60").Value = 39999
$ws.Range("K60").Value = 6498.4
$ws.Range("L60").Value = 39999
$ws.Range("M60").Value = -5987.4
$ws.Range("N60").Value = -41021

$ws.Range("H61").Value = 26665.5
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 29998.6
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 29998.6
$ws.Range("M61").Value = -9652
$ws.Range("N61").Value = -30694.6

$ws.Range("H62").Value = 7500
$ws.Range("J62").Value = 7500
$ws.Range("L62").Value = 7500
$ws.Range("N62").Value = -8748

$ws.Range("H65").Value = 7500
$ws.Range("J65").Value = 7500
$ws.Range("L65").Value = 37500
$ws.Range("N65").Value = -43740

$ws.Range("H105").Value = 8946.77
$ws.Range("I105").Value = 11030.8
$ws.Range("K105").Value = 11030.8
$ws.Range("M105").Value = -9283.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2471205
$ws.Range("I4").Value = 2728805.5
$ws.Range("K4").Value = 8186416.5
$ws.Range("M4").Value = -8186304.5

$ws.Range("H60").Value = 1509.8572
$ws.Range("I60").Value = 2085.8
$ws.Range("K60").Value = 6257.400000000001
$ws.Range("M60").Value = -6006.400000000001

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H113").Value = 8354.143
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 8354.143
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 25062.429
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -29402.429

$ws.Range("H120").Value = 8664.833000000001
$ws.Range("I120").Value = 5497.25
$ws.Range("J120").Value = 15000
$ws.Range("K120").Value = 16491.75
$ws.Range("L120").Value = 45000
$ws.Range("M120").Value = -11653.75
$ws.Range("N120").Value = -54676

$ws.Range("H131").Value = 1467.49
$ws.Range("I131").Value = 849.8
$ws.Range("J131").Value = 1500
$ws.Range("K131").Value = 2549.4
$ws.Range("L131").Value = 4500
$ws.Range("M131").Value = 2490.6
$ws.Range("N131").Value = -14580

$ws.Range("H140").Value = 3467.8333
$ws.Range("I140").Value = 1375.1428
$ws.Range("K140").Value = 4125.428400000001
$ws.Range("M140").Value = 1054.571599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 935778.5600000001
$ws.Range("I20").Value = 2506725
$ws.Range("J20").Value = 38094.855
$ws.Range("K20").Value = 2506725
$ws.Range("L20").Value = 38094.855
$ws.Range("M20").Value = -2506480
$ws.Range("N20").Value = -38584.855

$ws.Range("H24").Value = 9252239
$ws.Range("I24").Value = 5002000
$ws.Range("K24").Value = 5002000
$ws.Range("M24").Value = -5001827

$ws.Range("H113").Value = 29009
$ws.Range("I113").Value = 9011
$ws.Range("K113").Value = 9011
$ws.Range("M113").Value = -6841

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3322.6
$ws.Range("I46").Value = 700.5
$ws.Range("J46").Value = 3726
$ws.Range("K46").Value = 700.5
$ws.Range("L46").Value = 3726
$ws.Range("M46").Value = -512.5
$ws.Range("N46").Value = -4102

$ws.Range("H93").Value = 47624908
$ws.Range("I93").Value = 62507304
$ws.Range("J93").Value = 1234
$ws.Range("K93").Value = 62507304
$ws.Range("L93").Value = 1234
$ws.Range("M93").Value = -62506056
$ws.Range("N93").Value = -3730

$ws.Range("H122").Value = 24261400
$ws.Range("I122").Value = 35919196
$ws.Range("K122").Value = 107757588
$ws.Range("M122").Value = -107755138

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 12500
$ws.Range("I31").Value = 12500
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 12500
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -12152
$ws.Range("N31").ClearContents()

$ws.Range("H32").Value = 16499.5
$ws.Range("I32").Value = 15499.5
$ws.Range("J32").Value = 17499.5
$ws.Range("K32").Value = 15499.5
$ws.Range("L32").Value = 17499.5
$ws.Range("M32").Value = -15182.5
$ws.Range("N32").Value = -18133.5

$ws.Range("H34").Value = 192447
$ws.Range("J34").Value = 189929
$ws.Range("L34").Value = 189929
$ws.Range("N34").Value = -190335

$ws.Range("H37").Value = 48544.832
$ws.Range("I37").Value = 52793
$ws.Range("J37").Value = 46420.75
$ws.Range("K37").Value = 52793
$ws.Range("L37").Value = 46420.75
$ws.Range("M37").Value = -52590
$ws.Range("N37").Value = -46826.75

$ws.Range("H96").Value = 1039.4
$ws.Range("I96").Value = 899.5
$ws.Range("K96").Value = 899.5
$ws.Range("M96").Value = 473.5

$ws.Range("H116").Value = 170666.67
$ws.Range("J116").Value = 170666.67
$ws.Range("L116").Value = 170666.67
$ws.Range("N116").Value = -179844.67

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
